$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new 2024 data row at the top of the data table ---
# Inserting a row at row 2 pushes all existing year rows down by one,
# which matches every row in the sheet shifting from row N to row N+1.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = 2024
$ws.Range("B2").Value = 39
$ws.Range("C2").Value = 25
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 36
$ws.Range("F2").Value = 2956

# --- Correct the 2022 total (now on row 4 after the insert) ---
$ws.Range("F4").Value = 2964

# --- Misc view/formatting touch-ups made while updating the sheet ---
$ws.Columns.Item(1).ColumnWidth = 6.83203125

$win = $excel.ActiveWindow
$win.Zoom = 180
$ws.Range("G8").Select()

$ws.PageSetup.Orientation = 1
